# Applies the Sat Mar 16 11:29:08 UTC 2024 "cryptos list" refresh:
# updated prices / 1h volume % for (almost) every row, and three pairs of
# adjacent rows (10/11, 16/17, 48/49) swapped places in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.077.92"
$ws.Range("E2").Value = "  +1.69%  "

$ws.Range("D3").Value = "3.676.51"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.93"
$ws.Range("E5").Value = "  +6.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.52"
$ws.Range("E6").Value = "  +15.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +3.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.719"
$ws.Range("E9").Value = "  +4.67%  "

# Row 10: "Dogecoin" -> "Avalanche"
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.36"
$ws.Range("E10").Value = "  +19.94%  "

# Row 11: "Avalanche" -> "Dogecoin"
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.158"
$ws.Range("E11").Value = "  +1.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000284"
$ws.Range("E12").Value = "  +2.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.36"
$ws.Range("E13").Value = "  +1.93%  "

$ws.Range("D14").Value = "4.275.75"
$ws.Range("E14").Value = "  +1.35%  "

$ws.Range("D15").Value = "3.687.67"
$ws.Range("E15").Value = "  +0.68%  "

# Row 16: "Chainlink" -> "TRON"
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.127"
$ws.Range("E16").Value = "  +1.35%  "

# Row 17: "TRON" -> "Chainlink"
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.32"
$ws.Range("E17").Value = "  +2.08%  "

$ws.Range("E18").Value = "  +4.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  +1.65%  "

$ws.Range("D20").Value = "68.030.09"
$ws.Range("E20").Value = "  +1.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "405.68"
$ws.Range("E21").Value = "  +2.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.57"
$ws.Range("E22").Value = "  +4.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.36"
$ws.Range("E23").Value = "  +3.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.44"
$ws.Range("E24").Value = "  +9.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.05"
$ws.Range("E25").Value = "  +3.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.93"
$ws.Range("E26").Value = "  +4.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.02"
$ws.Range("E27").Value = "  +0.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.74"
$ws.Range("E28").Value = "  +1.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").Value = "  +3.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.40"
$ws.Range("E30").Value = "  +2.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.61"
$ws.Range("E31").Value = "  +5.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "47.25"
$ws.Range("E32").Value = "  +12.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.62"
$ws.Range("E33").Value = "  +4.44%  "

$ws.Range("E34").Value = "  +7.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "628.79"
$ws.Range("E35").Value = "  +8.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.03"
$ws.Range("E36").Value = "  +4.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.409"
$ws.Range("E37").Value = "  +7.00%  "

$ws.Range("D38").Value = "0.0₃0817"
$ws.Range("E38").Value = "  -5.38%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("E40").Value = "  -0.22%  "

$ws.Range("E41").Value = "  +5.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("E42").Value = "  +4.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0441"
$ws.Range("E43").Value = "  +4.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.61"
$ws.Range("E44").Value = "  +2.05%  "

$ws.Range("D45").Value = "2.904.86"
$ws.Range("E45").Value = "  +6.29%  "

$ws.Range("E46").Value = "  +6.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.19"
$ws.Range("E47").Value = "  +3.06%  "

# Row 48: "dogwifhat" -> "Monero"
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "146.12"
$ws.Range("E48").Value = "  +4.47%  "

# Row 49: "Monero" -> "dogwifhat"
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.68"
$ws.Range("E49").Value = "  -2.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.67"
$ws.Range("E50").Value = "  +2.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.06"
$ws.Range("E51").Value = "  -1.62%  "
